$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Rules" sheet rule #4 (row 11) used to be labelled with the text "R40".
# It is now re-labelled as the literal text "1". A leading apostrophe forces
# Excel to store the digit string as text (shared string) rather than as a
# number, matching the original cell's content type.
$ws.Range("B11").Value = "'1"
